$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated masking logic: refresh masked Name, Birth Date, and Home Address values
# for each student row (2-31). Birth Date moves from a styled numeric date value
# to a plain text (inline string) value, so we briefly mark the cell as Text before
# assigning it (to stop Excel re-parsing the string back into a date serial), then
# clear the temporary formatting so the cell keeps the default (unstyled) cell format.

$ws.Range("A2").Value = "Jennifer Rodriguez"
$dateCell = $ws.Range("C2")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1975-09-27"
$dateCell.ClearFormats()
$ws.Range("D2").Value = "92966 Brandon Canyon Suite 730`nDayland, IN 71365"

$ws.Range("A3").Value = "Susan Chapman"
$dateCell = $ws.Range("C3")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1972-02-12"
$dateCell.ClearFormats()
$ws.Range("D3").Value = "92484 Sanchez Burg Suite 331`nGrayland, OK 72314"

$ws.Range("A4").Value = "Tina Elliott"
$dateCell = $ws.Range("C4")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1926-05-04"
$dateCell.ClearFormats()
$ws.Range("D4").Value = "9373 Stacy Unions Suite 449`nCannonland, IA 21742"

$ws.Range("A5").Value = "Nicholas Wells"
$dateCell = $ws.Range("C5")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1925-11-10"
$dateCell.ClearFormats()
$ws.Range("D5").Value = "677 Ronald Lights`nAndrewchester, AS 30654"

$ws.Range("A6").Value = "Shawn Hunt"
$dateCell = $ws.Range("C6")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1976-07-02"
$dateCell.ClearFormats()
$ws.Range("D6").Value = "471 Bryan Fields`nGeorgestad, NJ 77162"

$ws.Range("A7").Value = "Kristen Fox"
$dateCell = $ws.Range("C7")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1957-01-29"
$dateCell.ClearFormats()
$ws.Range("D7").Value = "075 Macias Curve Suite 296`nLopezborough, NH 38025"

$ws.Range("A8").Value = "Madison Duncan"
$dateCell = $ws.Range("C8")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1947-02-06"
$dateCell.ClearFormats()
$ws.Range("D8").Value = "75232 Walters Pines`nJasonborough, OK 27600"

$ws.Range("A9").Value = "Laurie Holmes"
$dateCell = $ws.Range("C9")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1928-07-24"
$dateCell.ClearFormats()
$ws.Range("D9").Value = "595 Brenda Harbor Suite 060`nSouth Madisonfurt, AZ 49540"

$ws.Range("A10").Value = "Christina Espinoza"
$dateCell = $ws.Range("C10")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1975-03-10"
$dateCell.ClearFormats()
$ws.Range("D10").Value = "PSC 8157, Box 1769`nAPO AA 76606"

$ws.Range("A11").Value = "Brenda Petersen"
$dateCell = $ws.Range("C11")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2004-02-25"
$dateCell.ClearFormats()
$ws.Range("D11").Value = "566 Kevin Cliffs`nPort Carlos, AL 04739"

$ws.Range("A12").Value = "Sharon Owen"
$dateCell = $ws.Range("C12")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1939-12-02"
$dateCell.ClearFormats()
$ws.Range("D12").Value = "247 Felicia Keys`nPort Mary, NY 16973"

$ws.Range("A13").Value = "Ashley Gonzales"
$dateCell = $ws.Range("C13")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1979-02-24"
$dateCell.ClearFormats()
$ws.Range("D13").Value = "158 David Trafficway`nKatherinechester, AS 45943"

$ws.Range("A14").Value = "David Myers"
$dateCell = $ws.Range("C14")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1993-02-25"
$dateCell.ClearFormats()
$ws.Range("D14").Value = "4817 Tonya Meadows`nJohnstad, OK 61309"

$ws.Range("A15").Value = "Ricky Sherman"
$dateCell = $ws.Range("C15")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1957-04-19"
$dateCell.ClearFormats()
$ws.Range("D15").Value = "71763 Harris Point`nEast Casey, RI 14476"

$ws.Range("A16").Value = "Nicholas Petersen"
$dateCell = $ws.Range("C16")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1984-05-13"
$dateCell.ClearFormats()
$ws.Range("D16").Value = "381 Destiny Track`nEdwinview, ID 24102"

$ws.Range("A17").Value = "Daniel Garcia"
$dateCell = $ws.Range("C17")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1943-05-17"
$dateCell.ClearFormats()
$ws.Range("D17").Value = "PSC 3371, Box 6888`nAPO AE 68523"

$ws.Range("A18").Value = "Vicki Morgan"
$dateCell = $ws.Range("C18")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1983-09-16"
$dateCell.ClearFormats()
$ws.Range("D18").Value = "091 Fischer Dale`nAndreamouth, CT 33331"

$ws.Range("A19").Value = "Nicholas Davila"
$dateCell = $ws.Range("C19")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1939-09-09"
$dateCell.ClearFormats()
$ws.Range("D19").Value = "235 White Hills Apt. 307`nAndrewsville, GA 56368"

$ws.Range("A20").Value = "Alyssa Campos"
$dateCell = $ws.Range("C20")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1998-02-06"
$dateCell.ClearFormats()
$ws.Range("D20").Value = "384 Matthew Island`nTorresland, MO 57458"

$ws.Range("A21").Value = "James Jones"
$dateCell = $ws.Range("C21")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1976-05-01"
$dateCell.ClearFormats()
$ws.Range("D21").Value = "Unit 4917 Box 9699`nDPO AA 69452"

$ws.Range("A22").Value = "Jamie Adkins"
$dateCell = $ws.Range("C22")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1960-07-12"
$dateCell.ClearFormats()
$ws.Range("D22").Value = "0269 Timothy Island Suite 786`nNorth Joshua, NJ 16285"

$ws.Range("A23").Value = "Alyssa Torres"
$dateCell = $ws.Range("C23")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1995-12-22"
$dateCell.ClearFormats()
$ws.Range("D23").Value = "29927 Timothy Trafficway Apt. 293`nNorth Carla, HI 04674"

$ws.Range("A24").Value = "Dr. Ashley Padilla"
$dateCell = $ws.Range("C24")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1934-10-22"
$dateCell.ClearFormats()
$ws.Range("D24").Value = "90193 Diana Creek`nLeestad, DE 39643"

$ws.Range("A25").Value = "David Burch"
$dateCell = $ws.Range("C25")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1979-05-16"
$dateCell.ClearFormats()
$ws.Range("D25").Value = "44613 Jamie Orchard Apt. 408`nSparksstad, KY 54863"

$ws.Range("A26").Value = "Tiffany Riley"
$dateCell = $ws.Range("C26")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1995-09-17"
$dateCell.ClearFormats()
$ws.Range("D26").Value = "570 Rogers Locks Apt. 277`nEast Juliastad, PR 81704"

$ws.Range("A27").Value = "Curtis Joyce MD"
$dateCell = $ws.Range("C27")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1974-07-23"
$dateCell.ClearFormats()
$ws.Range("D27").Value = "PSC 3796, Box 8239`nAPO AA 82046"

$ws.Range("A28").Value = "Jessica Miller"
$dateCell = $ws.Range("C28")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1933-09-27"
$dateCell.ClearFormats()
$ws.Range("D28").Value = "4125 Bell Trail`nRyanfurt, AL 44201"

$ws.Range("A29").Value = "Alexa Ross"
$dateCell = $ws.Range("C29")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1946-07-30"
$dateCell.ClearFormats()
$ws.Range("D29").Value = "18124 Gary Cove Suite 762`nNew Jessicaport, CA 14465"

$ws.Range("A30").Value = "Mary Kelly"
$dateCell = $ws.Range("C30")
$dateCell.NumberFormat = "@"
$dateCell.Value = "1930-06-05"
$dateCell.ClearFormats()
$ws.Range("D30").Value = "USS Beard`nFPO AP 77606"

$ws.Range("A31").Value = "Mary Lopez"
$dateCell = $ws.Range("C31")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2001-03-16"
$dateCell.ClearFormats()
$ws.Range("D31").Value = "377 Melissa Cliff`nNew Kevin, NM 46953"
